# Generate Report for Handback
# Updates the Overview / zh-cn / de-de sheets of the localization-status report
# to reflect a failed handback transform for the 8f014f51-... file.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# 1. Update the "Ready for handoff" status text (shared string) to "Handback transform failed"
#    wherever it appears: Overview!E3, Overview!F3, zh-cn!C3, de-de!C3
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# 2. Populate the Error Detail column (P) for the 8f014f51 row (row 3) on zh-cn and de-de sheets
$zhcn.Range("P3").Value = "Handback file name: 0mnn51g3.340 is different with handoff file name: 8f014f51-3e08-47db-b122-55822e0e113d.3e34f6a30032a0cc1b66d733a3ab89431eb93285.zh-cn."
$dede.Range("P3").Value = "Handback file name: 0mnn51g3.340 is different with handoff file name: 8f014f51-3e08-47db-b122-55822e0e113d.3e34f6a30032a0cc1b66d733a3ab89431eb93285.de-de."

# 3. Widen the Error Detail column (P / column 16) on both language sheets so the new text is readable
#    (39.1667 in COM "characters" units renders as the target 40 "characters" in the saved OOXML
#    column width, which uses a slightly different pixel-rounded unit)
$zhcn.Columns.Item(16).ColumnWidth = 39.1667
$dede.Columns.Item(16).ColumnWidth = 39.1667
